$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 ("Work to Be Done"): merge "Each team member: f" + "ork " into a
# single run "Each team member: fork "
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(1).TextFrame.TextRange
$para10_1 = $tr10.Paragraphs(1,1)
$r10_1 = $para10_1.Runs(1,1)
$r10_2 = $para10_1.Runs(2,1)
$r10_1.Text = "Each team member: fork "
$r10_2.Text = ""

# ---------------------------------------------------------------------------
# Slide 2 ("Project Committee" -> "Project Team"): content placeholder edits
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Split "Nick Sellen - Technical Guru" run into " - Technical " + "Guru"
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$para2_3 = $tr2.Paragraphs(3,1)
$r2_3 = $para2_3.Runs(3,1)
$r2_3.Text = " – Technical "
$r2_3.InsertAfter("Guru") | Out-Null

# Add a new paragraph "You - Code Warriors!!!" after that paragraph
$tr2b = $s2.Shapes.Item(1).TextFrame.TextRange
$para2_3b = $tr2b.Paragraphs(3,1)
$para2_3b.InsertAfter("`rYou – ") | Out-Null

$tr2c = $s2.Shapes.Item(1).TextFrame.TextRange
$para2_4 = $tr2c.Paragraphs(4,1)
$r2_4_1 = $para2_4.Runs(1,1)
$r2_4_1.InsertAfter("Code Warriors!!!") | Out-Null

# Title: "Project Committee" -> "Project " + "Team"
$tr2title = $s2.Shapes.Item(2).TextFrame.TextRange
$para2title = $tr2title.Paragraphs(1,1)
$r2title = $para2title.Runs(1,1)
$r2title.Text = "Project "
$r2title.InsertAfter("Team") | Out-Null

# ---------------------------------------------------------------------------
# Slide 8 ("Project at a Glance (cont.)"): constraints wording + split bullet
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

# "Constraints:" -> "Constraints/Limitations:"
$tr8 = $s8.Shapes.Item(1).TextFrame.TextRange
$para8_1 = $tr8.Paragraphs(1,1)
$r8_1 = $para8_1.Runs(1,1)
$r8_1.Text = "Constraints/Limitations:"

# Split "Only worked in Firefox on my end (in web server, ran fine)" into two
# separate bullet paragraphs
$tr8b = $s8.Shapes.Item(1).TextFrame.TextRange
$para8_5 = $tr8b.Paragraphs(5,1)
$r8_5 = $para8_5.Runs(1,1)
$r8_5.Text = "Only worked in Firefox on my end "
$r8_5.InsertAfter("when not running web server") | Out-Null

$tr8c = $s8.Shapes.Item(1).TextFrame.TextRange
$para8_5b = $tr8c.Paragraphs(5,1)
$para8_5b.InsertAfter("`rRan fine on my end when running web server") | Out-Null

# ---------------------------------------------------------------------------
# Slide 9 ("Work Done So far"): merge the "Prepared this PowerPoint..." runs
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(1).TextFrame.TextRange
$para9_7 = $tr9.Paragraphs(7,1)
$r9_1 = $para9_7.Runs(1,1)
$r9_2 = $para9_7.Runs(2,1)
$r9_1.Text = "Prepared this PowerPoint project summary preview"
$r9_2.Text = ""
